# Applies the OOXML diff:
#  - Slide 10 ("Data mapping"): shift the text box + picture down (y only)
#  - Slide 6: shift the two big text boxes, and move the intro text box
#             ("Firstly started with a Student example ...") here from slide 8
#  - Slide 8 ("Student Example" -> "Requirement Explanation"): retitle,
#             shift the whole diagram up, re-center 3 labels, and drop the
#             intro text box (now living on slide 6)
#
# NOTE: Shape.Top/.Left are stored as single-precision (float32) points in
# this COM host (matching real PowerPoint VBA `Single` semantics), and EMU
# is recovered as floor(pt_as_float32 * 12700). The literals below were
# chosen so that round-trip lands exactly on the target EMU value from the
# reference diff (each is annotated with the EMU it reproduces).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10 - "Data mapping"
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)

$shp = $s10.Shapes.Item(2)                       # textbox, id 41
$shp.Top = 111.78126154251969                    # y -> 1419622 EMU

$shp = $s10.Shapes.Item(4)                       # picture, id 1026
$shp.Top = 107.9284251968504                     # y -> 1370691 EMU

# ---------------------------------------------------------------------
# Slide 8 - "Student Example" -> "Requirement Explanation"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Requirement Explanation"

$s8.Shapes.Item(2).Top = 173.52410128818897       # id 4  -> 2203756 EMU
$s8.Shapes.Item(3).Top = 173.48228346456693       # id 5  -> 2203225 EMU
$s8.Shapes.Item(4).Top = 173.48212598425198       # id 6  -> 2203223 EMU
$s8.Shapes.Item(5).Top = 268.6520538440945        # id 7  -> 3411881 EMU
$s8.Shapes.Item(6).Top = 266.9136220472441        # id 8  -> 3389803 EMU

# Move the intro text box (id 10) over to slide 6, then delete it here.
$introShape = $s8.Shapes.Item(8)
$introShape.Copy()
$s6 = $p.Slides.Item(6)
$pastedRange = $s6.Shapes.Paste()
$introCopy = $pastedRange.Item(1)
$introCopy.Left = 59.45244094488189               # x -> 755046 EMU
$introCopy.Top = 107.7663779527559                # y -> 1368633 EMU

$s8.Shapes.Item(9).Top = 231.304100088189         # id 24 connector -> 2937562 EMU
$s8.Shapes.Item(10).Top = 198.72181102362205      # id 35 connector -> 2523767 EMU
$s8.Shapes.Item(11).Top = 198.19496062992127      # id 40 connector -> 2517076 EMU
$s8.Shapes.Item(12).Top = 233.04267716535432      # id 54 connector -> 2959642 EMU

$s8.Shapes.Item(13).Top = 181.1852036503937       # id 61 "I"   -> 2301052 EMU
$s8.Shapes.Item(14).Top = 274.57464566929133      # id 62 "II"  -> 3487098 EMU
$s8.Shapes.Item(15).Top = 180.27291338582677      # id 63 "III" -> 2289466 EMU
$s8.Shapes.Item(16).Top = 273.07858267716534      # id 64 "IV"  -> 3468098 EMU
$s8.Shapes.Item(17).Top = 117.90378192755905      # id 65 "Linear data" -> 1497378 EMU

$shp66 = $s8.Shapes.Item(18)                      # id 66 "Linear data with unknown values"
$shp66.Left = 161.55275730551182                  # x -> 2051720 EMU
$shp66.Top = 322.02095038188975                   # y -> 4089666 EMU
$shp66.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$s8.Shapes.Item(19).Top = 117.90370078740158      # id 68 "2D data" -> 1497377 EMU

$shp69 = $s8.Shapes.Item(20)                      # id 69 "2D data with unknown values"
$shp69.Left = 371.339842519685                    # x -> 4716016 EMU
$shp69.Top = 322.54267886535433                   # y -> 4096292 EMU
$shp69.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$shp70 = $s8.Shapes.Item(21)                      # id 70 "Final requirements for the project"
$shp70.Left = 500.1828461456693                   # x -> 6352322 EMU
$shp70.Top = 106.56386186771654                   # y -> 1353361 EMU
$shp70.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$s8.Shapes.Item(22).Top = 180.79464566929133      # id 81 "V" -> 2296092 EMU

# Remove the intro text box from slide 8 now that it has been re-created on slide 6.
$introShape.Delete()

# ---------------------------------------------------------------------
# Slide 6
# ---------------------------------------------------------------------
$s6.Shapes.Item(3).Left = 76.50394060787401       # id 27 -> x 971600 EMU
$s6.Shapes.Item(3).Top = 151.4707107614173        # id 27 -> y 1923678 EMU

$s6.Shapes.Item(4).Left = 342.73292548582674      # id 29 -> x 4352708 EMU
$s6.Shapes.Item(4).Top = 151.4707107614173        # id 29 -> y 1923678 EMU
